# Insert a new "ListBullet" paragraph right after the
# "Docente(s) Responsável(eis)" heading, listing the four responsible
# lecturers, each on its own line (separated by manual line breaks),
# matching the upstream diff.

$d = $word.ActiveDocument

$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Docente(s) Responsável(eis)*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the 'Docente(s) Responsável(eis)' paragraph"
}

# Create the (initially empty) destination paragraph right after the
# heading paragraph.
$target.Range.InsertParagraphAfter() | Out-Null
$newPara = $target.Next()
$newRange = $newPara.Range

# Build the four names as their own runs (text + manual line break,
# except for the last entry) and apply the ListBullet paragraph style,
# all via a WordOpenXML fragment inserted into the fresh paragraph.
$ooxml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512">
<pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/></Relationships></pkg:xmlData>
</pkg:part>
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>
<w:p><w:pPr><w:pStyle w:val="ListBullet"/></w:pPr><w:r><w:t>471420 - Carlos Antonio Reis Pereira Baptista</w:t><w:br/></w:r><w:r><w:t>3480026 - Jo&#227;o Paulo Pascon</w:t><w:br/></w:r><w:r><w:t>5840793 - S&#233;rgio Schneider</w:t><w:br/></w:r><w:r><w:t>7797767 - Viktor Pastoukhov</w:t></w:r></w:p>
<w:sectPr/>
</w:body></w:document></pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$newRange.InsertXML($ooxml) | Out-Null
